# Update the division-fact table with newly generated values.
# Uses direct cell addressing (rather than text Find/Replace) because one of
# the source strings ("84÷7=12, 0") occurs twice in the table but must be
# replaced with two different results.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "13÷2=6, 1"
$t.Cell(1,2).Range.Text  = "36÷7=5, 1"
$t.Cell(1,3).Range.Text  = "91÷5=18, 1"
$t.Cell(1,4).Range.Text  = "77÷4=19, 1"
$t.Cell(1,5).Range.Text  = "33÷8=4, 1"

$t.Cell(5,1).Range.Text  = "21÷5=4, 1"
$t.Cell(5,2).Range.Text  = "19÷3=6, 1"
$t.Cell(5,3).Range.Text  = "56÷9=6, 2"
$t.Cell(5,4).Range.Text  = "78÷7=11, 1"
$t.Cell(5,5).Range.Text  = "73÷4=18, 1"

$t.Cell(9,1).Range.Text  = "80÷3=26, 2"
$t.Cell(9,2).Range.Text  = "75÷8=9, 3"
$t.Cell(9,3).Range.Text  = "78÷7=11, 1"
$t.Cell(9,4).Range.Text  = "88÷5=17, 3"
$t.Cell(9,5).Range.Text  = "52÷8=6, 4"

$t.Cell(13,1).Range.Text = "91÷3=30, 1"
$t.Cell(13,2).Range.Text = "70÷6=11, 4"
$t.Cell(13,3).Range.Text = "44÷4=11, 0"
$t.Cell(13,4).Range.Text = "70÷6=11, 4"
$t.Cell(13,5).Range.Text = "18÷8=2, 2"

$t.Cell(17,1).Range.Text = "98÷8=12, 2"
$t.Cell(17,2).Range.Text = "83÷7=11, 6"
$t.Cell(17,3).Range.Text = "80÷3=26, 2"
$t.Cell(17,4).Range.Text = "50÷3=16, 2"
$t.Cell(17,5).Range.Text = "58÷6=9, 4"
